# Append the 2025-12-26 data point to the "Chart" sheet and keep the
# "Table" sheet's header row pointing at the same (shared) strings.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

$usedRange = $chart.UsedRange
$priorLastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$lastRow = $priorLastRow + 1

# Work out the next calendar day from the last date already present in
# column A (stored as plain "yyyy-MM-dd" text).
$priorDateText = $chart.Range("A$priorLastRow").Text
$priorDate = [datetime]::ParseExact($priorDateText, "yyyy-MM-dd", $null)
$nextDateText = $priorDate.AddDays(1).ToString("yyyy-MM-dd")

# Write the new row's numeric values first.
$chart.Range("B$lastRow").Value = 0
$chart.Range("C$lastRow").Value = 32

# Write the date as literal text (not an auto-converted date serial) by
# using a leading apostrophe, matching how the existing date cells in
# column A are stored (shared string, general/default style).
$chart.Range("A$lastRow").Value = "'" + $nextDateText

# Re-apply the same cell format/style as the row above so the new cell
# keeps the default style (no quote-prefix style) like its neighbours.
$chart.Range("A$priorLastRow").Copy() | Out-Null
$chart.Range("A$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
